$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the existing "Landscape" data rows (22:26) down by one row (to 23:27) to
# make room for a new first row of data. Using Copy + PasteSpecial (rather than
# Rows.Insert) so the moved rows keep their per-row default style / customFormat
# metadata along with each cell's explicit style.
$ws.Range("A22:E26").Copy()
$ws.Range("A23:E27").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# The bottom-most row (27) is brand new to the sheet's dimensions, so make sure
# every one of its cells carries the same numeric-format style as column B.
$ws.Range("B27:E27").NumberFormat = "0.00000000000000000000"

# Populate the freshly-opened row 22 with its new values.
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").NumberFormat = "0.00000000000000000000"

$ws.Range("B22:E22").NumberFormat = "0.00000000000000000000"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 500
$ws.Range("D22").Value = -115
$ws.Range("E22").Value = 0.67895174917131795

# Match the author's final selection state.
$ws.Range("E18").Select()
